# Auto-generated: updates currentAveragePrice / Leve price & profit columns (H-N)
# across the 8 job sheets, per the scheduled pricing-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6333.3335
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 8500
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 8500
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -9748

$ws.Range("H65").Value = 6333.3335
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 8500
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 42500
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -48740

$ws.Range("H70").Value = 4225
$ws.Range("J70").Value = 4447.5
$ws.Range("L70").Value = 13342.5
$ws.Range("N70").Value = -13882.5

$ws.Range("H73").Value = 4225
$ws.Range("J73").Value = 4447.5
$ws.Range("L73").Value = 13342.5
$ws.Range("N73").Value = -15214.5

$ws.Range("H113").Value = 5929.154
$ws.Range("I113").Value = 4845
$ws.Range("J113").Value = 6126.273
$ws.Range("K113").Value = 4845
$ws.Range("L113").Value = 6126.273
$ws.Range("M113").Value = -1591
$ws.Range("N113").Value = -12634.273

$ws.Range("H137").Value = 2631.375
$ws.Range("I137").Value = 1846.0667
$ws.Range("J137").Value = 3940.2222
$ws.Range("K137").Value = 5538.2001
$ws.Range("L137").Value = 11820.6666
$ws.Range("M137").Value = -2988.2001
$ws.Range("N137").Value = -16920.6666


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1111.8889
$ws.Range("I45").Value = 1113.1428
$ws.Range("J45").Value = 1107.5
$ws.Range("K45").Value = 1113.1428
$ws.Range("L45").Value = 1107.5
$ws.Range("M45").Value = -736.1428000000001
$ws.Range("N45").Value = -1861.5

$ws.Range("H61").Value = 1793.4814
$ws.Range("I61").Value = 1817.35
$ws.Range("J61").Value = 1725.2858
$ws.Range("K61").Value = 1817.35
$ws.Range("L61").Value = 1725.2858
$ws.Range("M61").Value = -1605.35
$ws.Range("N61").Value = -2149.2858

$ws.Range("H74").Value = 8993.583000000001
$ws.Range("I74").Value = 13401.833
$ws.Range("J74").Value = 4585.3335
$ws.Range("K74").Value = 13401.833
$ws.Range("L74").Value = 4585.3335
$ws.Range("M74").Value = -12527.833
$ws.Range("N74").Value = -6333.3335

$ws.Range("H77").Value = 8993.583000000001
$ws.Range("I77").Value = 13401.833
$ws.Range("J77").Value = 4585.3335
$ws.Range("K77").Value = 67009.16500000001
$ws.Range("L77").Value = 22926.6675
$ws.Range("M77").Value = -62641.16500000001
$ws.Range("N77").Value = -31662.6675

$ws.Range("H80").Value = 45432.855
$ws.Range("J80").Value = 45432.855
$ws.Range("L80").Value = 45432.855
$ws.Range("N80").Value = -47428.855

$ws.Range("H83").Value = 45432.855
$ws.Range("J83").Value = 45432.855
$ws.Range("L83").Value = 136298.565
$ws.Range("N83").Value = -146282.565

$ws.Range("H136").Value = 1793.4814
$ws.Range("I136").Value = 1817.35
$ws.Range("J136").Value = 1725.2858
$ws.Range("K136").Value = 5452.049999999999
$ws.Range("L136").Value = 5175.857400000001
$ws.Range("M136").Value = -2902.049999999999
$ws.Range("N136").Value = -10275.8574


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 36739.25
$ws.Range("J82").Value = 46566.668
$ws.Range("L82").Value = 46566.668
$ws.Range("N82").Value = -47332.668

$ws.Range("H85").Value = 36739.25
$ws.Range("J85").Value = 46566.668
$ws.Range("L85").Value = 46566.668
$ws.Range("N85").Value = -49218.668

$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16669455
$ws.Range("I31").Value = 1431.9412
$ws.Range("J31").Value = 38466100
$ws.Range("K31").Value = 1431.9412
$ws.Range("L31").Value = 38466100
$ws.Range("M31").Value = -1136.9412
$ws.Range("N31").Value = -38466690

$ws.Range("H34").Value = 16669455
$ws.Range("I34").Value = 1431.9412
$ws.Range("J34").Value = 38466100
$ws.Range("K34").Value = 1431.9412
$ws.Range("L34").Value = 38466100
$ws.Range("M34").Value = -1229.9412
$ws.Range("N34").Value = -38466504

$ws.Range("H93").Value = 8000
$ws.Range("I93").Value = 8000
$ws.Range("K93").Value = 8000
$ws.Range("M93").Value = -6128

$ws.Range("H123").Value = 40733.332
$ws.Range("J123").Value = 40733.332
$ws.Range("L123").Value = 40733.332
$ws.Range("N123").Value = -50533.332

$ws.Range("H132").Value = 3065.2632
$ws.Range("I132").Value = 2044.2
$ws.Range("J132").Value = 4199.778
$ws.Range("K132").Value = 6132.6
$ws.Range("L132").Value = 12599.334
$ws.Range("M132").Value = -3602.6
$ws.Range("N132").Value = -17659.334


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 800
$ws.Range("J127").Value = 800
$ws.Range("L127").Value = 2400
$ws.Range("N127").Value = -12320

$ws.Range("H133").Value = 3334.4443
$ws.Range("I133").Value = 2905
$ws.Range("J133").Value = 4193.3335
$ws.Range("K133").Value = 8715
$ws.Range("L133").Value = 12580.0005
$ws.Range("M133").Value = -3655
$ws.Range("N133").Value = -22700.0005


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2988.5881
$ws.Range("I122").Value = 2220.4
$ws.Range("J122").Value = 8750
$ws.Range("K122").Value = 6661.200000000001
$ws.Range("L122").Value = 26250
$ws.Range("M122").Value = -4211.200000000001
$ws.Range("N122").Value = -31150

$ws.Range("H132").Value = 2333.976
$ws.Range("I132").Value = 1372.6428
$ws.Range("J132").Value = 4256.643
$ws.Range("K132").Value = 4117.928400000001
$ws.Range("L132").Value = 12769.929
$ws.Range("M132").Value = -1587.928400000001
$ws.Range("N132").Value = -17829.929


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 746.64703
$ws.Range("I16").Value = 746.64703
$ws.Range("K16").Value = 746.64703
$ws.Range("M16").Value = -576.64703

$ws.Range("H122").Value = 3031.8125
$ws.Range("I122").Value = 1700.3636
$ws.Range("J122").Value = 5961
$ws.Range("K122").Value = 5101.0908
$ws.Range("L122").Value = 17883
$ws.Range("M122").Value = -2651.0908
$ws.Range("N122").Value = -22783

$ws.Range("H136").Value = 4610.737
$ws.Range("I136").Value = 1430.4
$ws.Range("J136").Value = 8144.4443
$ws.Range("K136").Value = 4291.200000000001
$ws.Range("L136").Value = 24433.3329
$ws.Range("M136").Value = -1741.200000000001
$ws.Range("N136").Value = -29533.3329


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 37993.332
$ws.Range("I62").Value = 6990
$ws.Range("K62").Value = 6990
$ws.Range("M62").Value = -6366

$ws.Range("H65").Value = 37993.332
$ws.Range("I65").Value = 6990
$ws.Range("K65").Value = 34950
$ws.Range("M65").Value = -31830

